$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two informational values that used to live in A13/A14 into H2/H3 instead
$ws.Range("H2").Value = "Population size: 50"
$ws.Range("H3").Value = "Mutation Probability: 0.3"

# Remove the old standalone rows 13 and 14 (their content moved to H2/H3)
$ws.Range("A13").Value = $null
$ws.Range("A14").Value = $null
$ws.Rows("13:14").Delete()

# Update the selected cell shown when the workbook is opened
$ws.Range("A12").Select()
